$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.973.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.566.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.564.70"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.172.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.571.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.10"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.573.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.07"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.86%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.579"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.711.71"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.18"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.51%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.75%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.77%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +23.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.563.53"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.02"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.95"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.30"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.01"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0807"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +10.60%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.70"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.47"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.28%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.495.13"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.91"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.39"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +11.44%  "
